$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: remove now-unused columns U:AD from header rows 1 and 2 ---
$ws.Range("U1:AD2").Clear()

# --- Step 2: insert 4 new rows for Holden2.5/5/10/15 before the HexGrid rows ---
$ws.Rows("16:19").Insert()

# --- Step 3: fix up A/B columns for the (now shifted) rows 16-23 ---
$ws.Range("A16").Value = 14
$ws.Range("A16").Borders.LineStyle = 1
$ws.Range("B16").Value = "Holden2.5"
$ws.Range("A17").Value = 15
$ws.Range("A17").Borders.LineStyle = 1
$ws.Range("B17").Value = "Holden5"
$ws.Range("A18").Value = 16
$ws.Range("A18").Borders.LineStyle = 1
$ws.Range("B18").Value = "Holden10"
$ws.Range("A19").Value = 17
$ws.Range("A19").Borders.LineStyle = 1
$ws.Range("B19").Value = "Holden15"
$ws.Range("A20").Value = 18
$ws.Range("A20").Borders.LineStyle = 1
$ws.Range("B20").Value = "HexGrid-90degTilt2.5degRes"
$ws.Range("A21").Value = 19
$ws.Range("A21").Borders.LineStyle = 1
$ws.Range("B21").Value = "HexGrid-90degTilt5degRes"
$ws.Range("A22").Value = 20
$ws.Range("A22").Borders.LineStyle = 1
$ws.Range("B22").Value = "HexGrid-90degTilt10degRes"
$ws.Range("A23").Value = 21
$ws.Range("A23").Borders.LineStyle = 1
$ws.Range("B23").Value = "HexGrid-90degTilt15degRes"

# --- Step 4: write the updated numeric simulation data for rows 3-23, columns C:T ---
# row 3
$ws.Range("C3").Value = 0.9995050687056324
$ws.Range("D3").Value = 1.001974138940215
$ws.Range("E3").Value = 0.9995089791943357
$ws.Range("F3").Value = 1.001974138940215
$ws.Range("G3").Value = 0.9995089791943357
$ws.Range("H3").Value = 0.998679441927302
$ws.Range("I3").Value = 1.001086681666614
$ws.Range("J3").Value = 0.9995062097700599
$ws.Range("K3").Value = 0.9995089791943357
$ws.Range("L3").Value = 0.9995050687056324
$ws.Range("M3").Value = 1.000739603822923
$ws.Range("N3").Value = 1.000739603822923
$ws.Range("O3").Value = 1.000855296437487
$ws.Range("P3").Value = 1.000329395613394
$ws.Range("Q3").Value = 1.000329395613394
$ws.Range("R3").Value = 1.00012429150863
$ws.Range("S3").Value = 1.00012429150863
$ws.Range("T3").Value = 1.000043420034026

# row 4
$ws.Range("C4").Value = 0.9990445177703934
$ws.Range("D4").Value = 1.003811004089877
$ws.Range("E4").Value = 0.9990521657258106
$ws.Range("F4").Value = 1.003811004089877
$ws.Range("G4").Value = 0.9990521657258106
$ws.Range("H4").Value = 0.9974505937605196
$ws.Range("I4").Value = 1.002097821896416
$ws.Range("J4").Value = 0.9990467475825771
$ws.Range("K4").Value = 0.9990521657258106
$ws.Range("L4").Value = 0.9990445177703934
$ws.Range("M4").Value = 1.001427760930135
$ws.Range("N4").Value = 1.001427760930135
$ws.Range("O4").Value = 1.001651114585562
$ws.Range("P4").Value = 1.000635895862027
$ws.Range("Q4").Value = 1.000635895862027
$ws.Range("R4").Value = 1.000239963327973
$ws.Range("S4").Value = 1.000239963327973
$ws.Range("T4").Value = 1.000083808470932

# row 5
$ws.Range("C5").Value = 0.998169162602943
$ws.Range("D5").Value = 1.007303754946199
$ws.Range("E5").Value = 0.9981828813508146
$ws.Range("F5").Value = 1.007303754946199
$ws.Range("G5").Value = 0.9981828813508146
$ws.Range("H5").Value = 0.9951151586072533
$ws.Range("I5").Value = 1.00402024402196
$ws.Range("J5").Value = 0.9981731589032151
$ws.Range("K5").Value = 0.9981828813508146
$ws.Range("L5").Value = 0.998169162602943
$ws.Range("M5").Value = 1.002736458774571
$ws.Range("N5").Value = 1.002736458774571
$ws.Range("O5").Value = 1.003164387190367
$ws.Range("P5").Value = 1.001218599633319
$ws.Range("Q5").Value = 1.001218599633319
$ws.Range("R5").Value = 1.000459670062693
$ws.Range("S5").Value = 1.000459670062693
$ws.Range("T5").Value = 1.000160726738731

# row 6
$ws.Range("C6").Value = 0.9973148471175791
$ws.Range("D6").Value = 1.010712489876634
$ws.Range("E6").Value = 0.9973345171581478
$ws.Range("F6").Value = 1.010712489876634
$ws.Range("G6").Value = 0.9973345171581478
$ws.Range("H6").Value = 0.9928358624022283
$ws.Range("I6").Value = 1.00589641675649
$ws.Range("J6").Value = 0.9973205829757494
$ws.Range("K6").Value = 0.9973345171581478
$ws.Range("L6").Value = 0.9973148471175791
$ws.Range("M6").Value = 1.004013668497107
$ws.Range("N6").Value = 1.004013668497107
$ws.Range("O6").Value = 1.004641251250235
$ws.Range("P6").Value = 1.001787284717454
$ws.Range("Q6").Value = 1.001787284717454
$ws.Range("R6").Value = 1.000674092827627
$ws.Range("S6").Value = 1.000674092827627
$ws.Range("T6").Value = 1.000235786047805

# row 7
$ws.Range("C7").Value = 0.999968238490723
$ws.Range("D7").Value = 1.000127254257054
$ws.Range("E7").Value = 0.9999680930491397
$ws.Range("F7").Value = 1.000127254257054
$ws.Range("G7").Value = 0.9999680930491397
$ws.Range("H7").Value = 0.999915333972783
$ws.Range("I7").Value = 1.000069956659459
$ws.Range("J7").Value = 0.999968197093725
$ws.Range("K7").Value = 0.9999680930491397
$ws.Range("L7").Value = 0.999968238490723
$ws.Range("M7").Value = 1.000047746373889
$ws.Range("N7").Value = 1.000047746373889
$ws.Range("O7").Value = 1.000055149802412
$ws.Range("P7").Value = 1.000021195265639
$ws.Range("Q7").Value = 1.000021195265639
$ws.Range("R7").Value = 1.000007919711514
$ws.Range("S7").Value = 1.000007919711514
$ws.Range("T7").Value = 1.000002845587147

# row 8
$ws.Range("C8").Value = 0.9999285719978939
$ws.Range("D8").Value = 1.00028575229939
$ws.Range("E8").Value = 0.999928544456024
$ws.Range("F8").Value = 1.00028575229939
$ws.Range("G8").Value = 0.999928544456024
$ws.Range("H8").Value = 0.9998095337126099
$ws.Range("I8").Value = 1.000157157804681
$ws.Range("J8").Value = 0.9999285646702367
$ws.Range("K8").Value = 0.999928544456024
$ws.Range("L8").Value = 0.9999285719978939
$ws.Range("M8").Value = 1.000107162148642
$ws.Range("N8").Value = 1.000107162148642
$ws.Range("O8").Value = 1.000123827367321
$ws.Range("P8").Value = 1.000047622917769
$ws.Range("Q8").Value = 1.000047622917769
$ws.Range("R8").Value = 1.000017853302333
$ws.Range("S8").Value = 1.000017853302333
$ws.Range("T8").Value = 1.000006354156806

# row 9
$ws.Range("C9").Value = 0.9999049373744728
$ws.Range("D9").Value = 1.000379257348897
$ws.Range("E9").Value = 0.999905632835272
$ws.Range("F9").Value = 1.000379257348897
$ws.Range("G9").Value = 0.999905632835272
$ws.Range("H9").Value = 0.9997463698269369
$ws.Range("I9").Value = 1.000208753503142
$ws.Range("J9").Value = 0.9999051406293841
$ws.Range("K9").Value = 0.999905632835272
$ws.Range("L9").Value = 0.9999049373744728
$ws.Range("M9").Value = 1.000142097361685
$ws.Range("N9").Value = 1.000142097361685
$ws.Range("O9").Value = 1.000164316075504
$ws.Range("P9").Value = 1.000063275852881
$ws.Range("Q9").Value = 1.000063275852881
$ws.Range("R9").Value = 1.000023865098479
$ws.Range("S9").Value = 1.000023865098479
$ws.Range("T9").Value = 1.000008348586351

# row 10
$ws.Range("C10").Value = 0.9997900146754395
$ws.Range("D10").Value = 1.000839106437341
$ws.Range("E10").Value = 0.9997906000417703
$ws.Range("F10").Value = 1.000839106437341
$ws.Range("G10").Value = 0.9997906000417703
$ws.Range("H10").Value = 0.9994399286414379
$ws.Range("I10").Value = 1.000461643950858
$ws.Range("J10").Value = 0.9997901865837786
$ws.Range("K10").Value = 0.9997906000417703
$ws.Range("L10").Value = 0.9997900146754395
$ws.Range("M10").Value = 1.00031456055639
$ws.Range("N10").Value = 1.00031456055639
$ws.Range("O10").Value = 1.000363588354546
$ws.Range("P10").Value = 1.000139907051517
$ws.Range("Q10").Value = 1.000139907051517
$ws.Range("R10").Value = 1.00005258029908
$ws.Range("S10").Value = 1.00005258029908
$ws.Range("T10").Value = 1.000018580055104

# row 11
$ws.Range("C11").Value = 0.9996676153769876
$ws.Range("D11").Value = 1.001324736609869
$ws.Range("E11").Value = 0.9996709748725832
$ws.Range("F11").Value = 1.001324736609869
$ws.Range("G11").Value = 0.9996709748725832
$ws.Range("H11").Value = 0.999113005494895
$ws.Range("I11").Value = 1.000729384138211
$ws.Range("J11").Value = 0.9996685970698705
$ws.Range("K11").Value = 0.9996709748725832
$ws.Range("L11").Value = 0.9996676153769876
$ws.Range("M11").Value = 1.000496175993428
$ws.Range("N11").Value = 1.000496175993428
$ws.Range("O11").Value = 1.000573912041689
$ws.Range("P11").Value = 1.000221108953147
$ws.Range("Q11").Value = 1.000221108953147
$ws.Range("R11").Value = 1.000083575433006
$ws.Range("S11").Value = 1.000083575433006
$ws.Range("T11").Value = 1.000029052260403

# row 12
$ws.Range("C12").Value = 1.001628178870439
$ws.Range("D12").Value = 0.9930606403511493
$ws.Range("E12").Value = 1.001926814300016
$ws.Range("F12").Value = 0.9930606403511493
$ws.Range("G12").Value = 1.001926814300016
$ws.Range("H12").Value = 1.004284944662976
$ws.Range("I12").Value = 0.9962524615799423
$ws.Range("J12").Value = 1.001715254269505
$ws.Range("K12").Value = 1.001926814300016
$ws.Range("L12").Value = 1.001628178870439
$ws.Range("M12").Value = 0.997344409610794
$ws.Range("N12").Value = 0.997344409610794
$ws.Range("O12").Value = 0.9969804269338435
$ws.Range("P12").Value = 0.9988718778405344
$ws.Range("Q12").Value = 0.9988718778405344
$ws.Range("R12").Value = 0.9996356119554047
$ws.Range("S12").Value = 0.9996356119554047
$ws.Range("T12").Value = 0.9998113823390044

# row 13
$ws.Range("C13").Value = 1.000463792942587
$ws.Range("D13").Value = 0.9983240710954566
$ws.Range("E13").Value = 1.000338328286865
$ws.Range("F13").Value = 0.9983240710954566
$ws.Range("G13").Value = 1.000338328286865
$ws.Range("H13").Value = 1.001260677753924
$ws.Range("I13").Value = 0.9990492061484643
$ws.Range("J13").Value = 1.000427214557603
$ws.Range("K13").Value = 1.000338328286865
$ws.Range("L13").Value = 1.000463792942587
$ws.Range("M13").Value = 0.999393932019022
$ws.Range("N13").Value = 0.999393932019022
$ws.Range("O13").Value = 0.9992790233955028
$ws.Range("P13").Value = 0.9997087307749696
$ws.Range("Q13").Value = 0.9997087307749696
$ws.Range("R13").Value = 0.9998661301529435
$ws.Range("S13").Value = 0.9998661301529435
$ws.Range("T13").Value = 0.9999772151308166

# row 14
$ws.Range("C14").Value = 1.002724931667661
$ws.Range("D14").Value = 0.9897079552393377
$ws.Range("E14").Value = 1.00229955870053
$ws.Range("F14").Value = 0.9897079552393377
$ws.Range("G14").Value = 1.00229955870053
$ws.Range("H14").Value = 1.007347508454681
$ws.Range("I14").Value = 0.9942409311711703
$ws.Range("J14").Value = 1.002600916710594
$ws.Range("K14").Value = 1.00229955870053
$ws.Range("L14").Value = 1.002724931667661
$ws.Range("M14").Value = 0.9962164434534996
$ws.Range("N14").Value = 0.9962164434534996
$ws.Range("O14").Value = 0.9955579393593897
$ws.Range("P14").Value = 0.9982441485358432
$ws.Range("Q14").Value = 0.9982441485358432
$ws.Range("R14").Value = 0.9992580010770149
$ws.Range("S14").Value = 0.9992580010770149
$ws.Range("T14").Value = 0.9998203003239957

# row 15
$ws.Range("C15").Value = 1.000606737637396
$ws.Range("D15").Value = 0.9973229015443735
$ws.Range("E15").Value = 1.00078182513898
$ws.Range("F15").Value = 0.9973229015443735
$ws.Range("G15").Value = 1.00078182513898
$ws.Range("H15").Value = 1.001584642071828
$ws.Range("I15").Value = 0.9985681129572284
$ws.Range("J15").Value = 1.00065779050995
$ws.Range("K15").Value = 1.00078182513898
$ws.Range("L15").Value = 1.000606737637396
$ws.Range("M15").Value = 0.9989648195908847
$ws.Range("N15").Value = 0.9989648195908847
$ws.Range("O15").Value = 0.9988325840463327
$ws.Range("P15").Value = 0.9995704881069165
$ws.Range("Q15").Value = 0.9995704881069164
$ws.Range("R15").Value = 0.9998733223649323
$ws.Range("S15").Value = 0.9998733223649323
$ws.Range("T15").Value = 0.9999203349766259

# row 16
$ws.Range("C16").Value = 0.9903509073044839
$ws.Range("D16").Value = 1.038501574190725
$ws.Range("E16").Value = 0.9904172594852118
$ws.Range("F16").Value = 1.038501574190725
$ws.Range("G16").Value = 0.9904172594852118
$ws.Range("H16").Value = 0.9742564560791315
$ws.Range("I16").Value = 1.021191219788074
$ws.Range("J16").Value = 0.9903702525280895
$ws.Range("K16").Value = 0.9904172594852118
$ws.Range("L16").Value = 0.9903509073044839
$ws.Range("M16").Value = 1.014426240747605
$ws.Range("N16").Value = 1.014426240747605
$ws.Range("O16").Value = 1.016681233761094
$ws.Range("P16").Value = 1.006423246993474
$ws.Range("Q16").Value = 1.006423246993474
$ws.Range("R16").Value = 1.002421750116408
$ws.Range("S16").Value = 1.002421750116408
$ws.Range("T16").Value = 1.000847944895953

# row 17
$ws.Range("C17").Value = 0.99211599716858
$ws.Range("D17").Value = 1.031487071664842
$ws.Range("E17").Value = 0.9921502510957574
$ws.Range("F17").Value = 1.031487071664842
$ws.Range("G17").Value = 0.9921502510957574
$ws.Range("H17").Value = 0.9789694746705115
$ws.Range("I17").Value = 1.017325815791424
$ws.Range("J17").Value = 0.99212598439397
$ws.Range("K17").Value = 0.9921502510957574
$ws.Range("L17").Value = 0.99211599716858
$ws.Range("M17").Value = 1.011801534416711
$ws.Range("N17").Value = 1.011801534416711
$ws.Range("O17").Value = 1.013642961541615
$ws.Range("P17").Value = 1.00525110664306
$ws.Range("Q17").Value = 1.00525110664306
$ws.Range("R17").Value = 1.001975892756234
$ws.Range("S17").Value = 1.001975892756234
$ws.Range("T17").Value = 1.000695765797514

# row 18
$ws.Range("C18").Value = 0.9956780518044105
$ws.Range("D18").Value = 1.017328676008452
$ws.Range("E18").Value = 0.9956494357383376
$ws.Range("F18").Value = 1.017328676008452
$ws.Range("G18").Value = 0.9956494357383376
$ws.Range("H18").Value = 0.9884802527515555
$ws.Range("I18").Value = 1.009524150874993
$ws.Range("J18").Value = 0.9956697080117675
$ws.Range("K18").Value = 0.9956494357383376
$ws.Range("L18").Value = 0.9956780518044105
$ws.Range("M18").Value = 1.006503363906431
$ws.Range("N18").Value = 1.006503363906431
$ws.Range("O18").Value = 1.007510292895952
$ws.Range("P18").Value = 1.0028853878504
$ws.Range("Q18").Value = 1.0028853878504
$ws.Range("R18").Value = 1.001076399822384
$ws.Range("S18").Value = 1.001076399822384
$ws.Range("T18").Value = 1.000388379198253

# row 19
$ws.Range("C19").Value = 0.9953072836834721
$ws.Range("D19").Value = 1.018836570298831
$ws.Range("E19").Value = 0.995261293062637
$ws.Range("F19").Value = 1.018836570298831
$ws.Range("G19").Value = 0.995261293062637
$ws.Range("H19").Value = 0.9874948513908497
$ws.Range("I19").Value = 1.010349471875714
$ws.Range("J19").Value = 0.9952938737281763
$ws.Range("K19").Value = 0.995261293062637
$ws.Range("L19").Value = 0.9953072836834721
$ws.Range("M19").Value = 1.007071926991151
$ws.Range("N19").Value = 1.007071926991151
$ws.Range("O19").Value = 1.008164441952672
$ws.Range("P19").Value = 1.00313504901498
$ws.Range("Q19").Value = 1.00313504901498
$ws.Range("R19").Value = 1.001166610026894
$ws.Range("S19").Value = 1.001166610026894
$ws.Range("T19").Value = 1.00042389067328

# row 20
$ws.Range("C20").Value = 0.9999982201373447
$ws.Range("D20").Value = 1.000008667036453
$ws.Range("E20").Value = 0.999997137140615
$ws.Range("F20").Value = 1.000008667036453
$ws.Range("G20").Value = 0.999997137140615
$ws.Range("H20").Value = 0.9999954634901617
$ws.Range("I20").Value = 1.00000451677198
$ws.Range("J20").Value = 0.9999979053407558
$ws.Range("K20").Value = 0.999997137140615
$ws.Range("L20").Value = 0.9999982201373447
$ws.Range("M20").Value = 1.000003443586899
$ws.Range("N20").Value = 1.000003443586899
$ws.Range("O20").Value = 1.000003801315259
$ws.Range("P20").Value = 1.000001341438137
$ws.Range("Q20").Value = 1.000001341438137
$ws.Range("R20").Value = 1.000000290363757
$ws.Range("S20").Value = 1.000000290363757
$ws.Range("T20").Value = 1.000000318319552

# row 21
$ws.Range("C21").Value = 0.9999565484448533
$ws.Range("D21").Value = 1.000172314683491
$ws.Range("E21").Value = 0.9999575928523742
$ws.Range("F21").Value = 1.000172314683491
$ws.Range("G21").Value = 0.9999575928523742
$ws.Range("H21").Value = 0.999883933753003
$ws.Range("I21").Value = 1.000095015443496
$ws.Range("J21").Value = 0.9999568539160034
$ws.Range("K21").Value = 0.9999575928523742
$ws.Range("L21").Value = 0.9999565484448533
$ws.Range("M21").Value = 1.000064431564172
$ws.Range("N21").Value = 1.000064431564172
$ws.Range("O21").Value = 1.000074626190613
$ws.Range("P21").Value = 1.00002881866024
$ws.Range("Q21").Value = 1.00002881866024
$ws.Range("R21").Value = 1.000011012208273
$ws.Range("S21").Value = 1.000011012208273
$ws.Range("T21").Value = 1.00000370984887

# row 22
$ws.Range("C22").Value = 0.9998482435396185
$ws.Range("D22").Value = 1.000609039630626
$ws.Range("E22").Value = 0.9998468388062062
$ws.Range("F22").Value = 1.000609039630626
$ws.Range("G22").Value = 0.9998468388062062
$ws.Range("H22").Value = 0.9995955785206447
$ws.Range("I22").Value = 1.000334647094415
$ws.Range("J22").Value = 0.9998478323357067
$ws.Range("K22").Value = 0.9998468388062062
$ws.Range("L22").Value = 0.9998482435396185
$ws.Range("M22").Value = 1.000228641585122
$ws.Range("N22").Value = 1.000228641585122
$ws.Range("O22").Value = 1.000263976754886
$ws.Range("P22").Value = 1.00010137399215
$ws.Range("Q22").Value = 1.00010137399215
$ws.Range("R22").Value = 1.000037740195664
$ws.Range("S22").Value = 1.000037740195664
$ws.Range("T22").Value = 1.000013696654536

# row 23
$ws.Range("C23").Value = 0.9996450872022719
$ws.Range("D23").Value = 1.001427991190117
$ws.Range("E23").Value = 0.9996392566565201
$ws.Range("F23").Value = 1.001427991190117
$ws.Range("G23").Value = 0.9996392566565201
$ws.Range("H23").Value = 0.9990546828882355
$ws.Range("I23").Value = 1.000784040606523
$ws.Range("J23").Value = 0.9996433842722557
$ws.Range("K23").Value = 0.9996392566565201
$ws.Range("L23").Value = 0.9996450872022719
$ws.Range("M23").Value = 1.000536539196195
$ws.Range("N23").Value = 1.000536539196195
$ws.Range("O23").Value = 1.000619039666304
$ws.Range("P23").Value = 1.000237445016303
$ws.Range("Q23").Value = 1.000237445016303
$ws.Range("R23").Value = 1.000087897926357
$ws.Range("S23").Value = 1.000087897926357
$ws.Range("T23").Value = 1.000032407135987

